# Edit: insert a new weekly price record as row 37, shifting all existing
# records (previously rows 37-158) down by one row (to rows 38-159).
#
# This matches the diff: a new row is inserted in the middle of the data
# table (right before the former row 37), so every subsequent row's content
# is now the content of the row that used to precede it, and the very last
# row (159) ends up with the data that used to be in row 158. The sheet's
# dimension grows from A1:R158 to A1:R159.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 37; this pushes old rows 37..158 down
# to 38..159 (Excel shifts formatting/content automatically).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new data record.
$ws.Range("A37").Value2 = 4
$ws.Range("B37").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value2 = "Los Lagos"
$ws.Range("D37").Value2 = 44481
$ws.Range("E37").Value2 = 10
$ws.Range("F37").Value2 = 100112044
$ws.Range("G37").Value2 = "Perejil"
$ws.Range("H37").Value2 = "Sin especificar"
$ws.Range("I37").Value2 = "Primera"
$ws.Range("J37").Value2 = 180
$ws.Range("K37").Value2 = 4500
$ws.Range("L37").Value2 = 4500
$ws.Range("M37").Value2 = 4500
$ws.Range("N37").Value2 = "`$/docena de atados (3 kilos)"
$ws.Range("O37").Value2 = "Región Metropolitana"
$ws.Range("P37").Value2 = 1500
$ws.Range("Q37").Value2 = 3
$ws.Range("R37").Value2 = "Hortaliza"
